# Jean's suggestion on Human in the loop for consideration
#
# Adds two new survey respondents (rows 36-37) to Sheet1, matching the
# column widths used for the main data columns, and leaves the selection
# on the newly-added data (L28) with the view scrolled down to row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A, B, C, E) -------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.830729166666668
$ws.Columns.Item(2).ColumnWidth = 21.385416666666668
$ws.Columns.Item(3).ColumnWidth = 21.276041666666668
$ws.Columns.Item(5).ColumnWidth = 52.498697916666664

# --- Row 36: Dr. John Williams ---------------------------------------------
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(36,1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(36,1).Value = 45418.502638888902
$ws.Cells.Item(36,2).Value = "Dr. John Williams"
$ws.Cells.Item(36,3).Value = "Male"
$ws.Cells.Item(36,4).Value = "21-30"
$ws.Cells.Item(36,5).Value = "test@abc.com.sg"
$ws.Cells.Item(36,6).Value = 30338111

# --- Row 37: Candice Washington ---------------------------------------------
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(37,1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(37,1).Value = 45419.1320949074
$ws.Cells.Item(37,2).Value = "Candice Washington"
$ws.Cells.Item(37,3).Value = "Female"
$ws.Cells.Item(37,4).Value = "41-50"
$ws.Cells.Item(37,5).Value = "test@abc.com.sg"
$ws.Cells.Item(37,6).Value = 87857672

# --- View / selection -------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L28").Select() | Out-Null
